$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "/src/Excel/entregable2/Mod_Masiva_Clientes_PN_18.csv"
$ws.Range("E2").Value = "PASSED"
$ws.Range("F2").Value = "28 jun. 2023, 14:38:44"

$ws.Columns.Item(3).ColumnWidth = 48
$ws.Range("C8").Select() | Out-Null
